# Daily attendance processing - 2025-12-17 04:17:57
# Normalises the "Recorded By" column (G): when a cell lists both the
# "System" actor and a real recorder address (backup@backdoor.com or
# dnasr281@gmail.com), the comma-separated list is reordered so the
# real address comes first.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $text = $cell.Text

    if ($text -eq $null -or $text -eq "") {
        continue
    }

    $parts = $text -split ", "
    $count = $parts.Length

    if ($count -le 1) {
        continue
    }

    $hasSystem = $false
    $hasTarget = $false
    foreach ($part in $parts) {
        if ($part -eq "System") {
            $hasSystem = $true
        }
        if ($part -eq "backup@backdoor.com" -or $part -eq "dnasr281@gmail.com") {
            $hasTarget = $true
        }
    }

    if ($hasSystem -and $hasTarget) {
        $reversed = @()
        for ($i = $count - 1; $i -ge 0; $i--) {
            $reversed += $parts[$i]
        }
        $newValue = [string]::Join(", ", $reversed)
        $cell.Value = $newValue
    }
}
